# Weekly Fruit/Vegetable price update ("Fruta / hortaliza, semanal").
#
# A new weekly price observation is prepended to the historical series:
# row 53 (first data row of this particular market/product block) is
# pushed down, shifting every subsequent record down by one row, and a
# brand-new record is written into the freed-up row 53. The table grows
# from 137 to 138 data rows (dimension A1:R137 -> A1:R138).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 53:137 down to 54:138, leaving row 53 empty for the new record.
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new weekly observation.
$ws.Cells.Item(53, 1).Value = 11
$ws.Cells.Item(53, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(53, 3).Value = "Bíobío"
$ws.Cells.Item(53, 4).Value = 45174
$ws.Cells.Item(53, 5).Value = 8
$ws.Cells.Item(53, 6).Value = 100112012
$ws.Cells.Item(53, 7).Value = "Espinaca"
$ws.Cells.Item(53, 8).Value = "Sin especificar"
$ws.Cells.Item(53, 9).Value = "Primera"
$ws.Cells.Item(53, 10).Value = 50
$ws.Cells.Item(53, 11).Value = 7500
$ws.Cells.Item(53, 12).Value = 7500
$ws.Cells.Item(53, 13).Value = 7500
$ws.Cells.Item(53, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(53, 15).Value = "Región Metropolitana"
$ws.Cells.Item(53, 16).Value = 750
$ws.Cells.Item(53, 17).Value = 10
$ws.Cells.Item(53, 18).Value = "Hortaliza"
